$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.348.04"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "3.540.79"
$ws.Range("E3").Value = "  -2.76%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.74"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.95"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("D8").Value = "3.537.74"
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.583"
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.49"
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000277"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "4.110.94"
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.54"
$ws.Range("E16").Value = "  -3.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "629.96"
$ws.Range("E17").Value = "  -6.02%  "
$ws.Range("D18").Value = "3.541.86"
$ws.Range("E18").Value = "  -3.05%  "
$ws.Range("D19").Value = "69.452.02"
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.19"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.891"
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.97"
$ws.Range("E24").Value = "  -6.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.44"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.82"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -4.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("E29").Value = "  -6.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.90"
$ws.Range("E30").Value = "  -5.41%  "
$ws.Range("E31").Value = "  -3.98%  "
$ws.Range("E32").Value = "  -5.57%  "
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.00"
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "637.26"
$ws.Range("E35").Value = "  +9.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.81"
$ws.Range("E36").Value = "  -2.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.52"
$ws.Range("E37").Value = "  -11.25%  "
$ws.Range("E38").Value = "  -3.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.34"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0458"
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").Value = "3.392.57"
$ws.Range("E43").Value = "  -5.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.329"
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0703"
$ws.Range("E45").Value = "  -5.11%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.86"
$ws.Range("E46").Value = "  -5.58%  "
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("E48").Value = "  -5.79%  "
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.45"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("E51").Value = "  +13.91%  "
